$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 96.28570999999999  # H5: 100.5 -> 96.28570999999999
$ws.Cells.Item(5, 9).Value = 96.28570999999999  # I5: 100.5 -> 96.28570999999999
$ws.Cells.Item(5, 11).Value = 96.28570999999999  # K5: 100.5 -> 96.28570999999999
$ws.Cells.Item(5, 13).Value = 18.71429000000001  # M5: 14.5 -> 18.71429000000001
$ws.Cells.Item(9, 8).Value = 565.4286  # H9: 453.05554 -> 565.4286
$ws.Cells.Item(9, 9).Value = 455.07693  # I9: 362.05884 -> 455.07693
$ws.Cells.Item(9, 11).Value = 455.07693  # K9: 362.05884 -> 455.07693
$ws.Cells.Item(9, 13).Value = -286.07693  # M9: -193.05884 -> -286.07693
$ws.Cells.Item(55, 8).Value = 314.25  # H55: 364.45456 -> 314.25
$ws.Cells.Item(55, 9).Value = 183.4  # I55: 183.6 -> 183.4
$ws.Cells.Item(55, 10).Value = 532.3333  # J55: 515.1667 -> 532.3333
$ws.Cells.Item(55, 11).Value = 183.4  # K55: 183.6 -> 183.4
$ws.Cells.Item(55, 12).Value = 532.3333  # L55: 515.1667 -> 532.3333
$ws.Cells.Item(55, 13).Value = 30.59999999999999  # M55: 30.40000000000001 -> 30.59999999999999
$ws.Cells.Item(55, 14).Value = -960.3333  # N55: -943.1667 -> -960.3333
$ws.Cells.Item(80, 8).Value = 1937.25  # H80: 2024.6875 -> 1937.25
$ws.Cells.Item(80, 9).Value = 761  # I80: 958.6 -> 761
$ws.Cells.Item(80, 10).Value = 2329.3333  # J80: 2509.2727 -> 2329.3333
$ws.Cells.Item(80, 11).Value = 2283  # K80: 2875.8 -> 2283
$ws.Cells.Item(80, 12).Value = 6987.999899999999  # L80: 7527.8181 -> 6987.999899999999
$ws.Cells.Item(80, 13).Value = -1285  # M80: -1877.8 -> -1285
$ws.Cells.Item(80, 14).Value = -8983.999899999999  # N80: -9523.8181 -> -8983.999899999999
$ws.Cells.Item(83, 8).Value = 1937.25  # H83: 2024.6875 -> 1937.25
$ws.Cells.Item(83, 9).Value = 761  # I83: 958.6 -> 761
$ws.Cells.Item(83, 10).Value = 2329.3333  # J83: 2509.2727 -> 2329.3333
$ws.Cells.Item(83, 11).Value = 6849  # K83: 8627.4 -> 6849
$ws.Cells.Item(83, 12).Value = 20963.9997  # L83: 22583.4543 -> 20963.9997
$ws.Cells.Item(83, 13).Value = -1857  # M83: -3635.4 -> -1857
$ws.Cells.Item(83, 14).Value = -30947.9997  # N83: -32567.4543 -> -30947.9997
$ws.Cells.Item(86, 8).Value = 3473.9412  # H86: 3753 -> 3473.9412
$ws.Cells.Item(86, 9).Value = 3716.6667  # I86: 3841.25 -> 3716.6667
$ws.Cells.Item(86, 10).Value = 2891.4  # J86: 3488.25 -> 2891.4
$ws.Cells.Item(86, 11).Value = 3716.6667  # K86: 3841.25 -> 3716.6667
$ws.Cells.Item(86, 12).Value = 2891.4  # L86: 3488.25 -> 2891.4
$ws.Cells.Item(86, 13).Value = -2593.6667  # M86: -2718.25 -> -2593.6667
$ws.Cells.Item(86, 14).Value = -5137.4  # N86: -5734.25 -> -5137.4
$ws.Cells.Item(88, 8).Value = 9987.4  # H88: 9992 -> 9987.4
$ws.Cells.Item(88, 10).Value = 9987.4  # J88: 9992 -> 9987.4
$ws.Cells.Item(88, 12).Value = 9987.4  # L88: 9992 -> 9987.4
$ws.Cells.Item(88, 14).Value = -10799.4  # N88: -10804 -> -10799.4
$ws.Cells.Item(89, 8).Value = 3473.9412  # H89: 3753 -> 3473.9412
$ws.Cells.Item(89, 9).Value = 3716.6667  # I89: 3841.25 -> 3716.6667
$ws.Cells.Item(89, 10).Value = 2891.4  # J89: 3488.25 -> 2891.4
$ws.Cells.Item(89, 11).Value = 18583.3335  # K89: 19206.25 -> 18583.3335
$ws.Cells.Item(89, 12).Value = 14457  # L89: 17441.25 -> 14457
$ws.Cells.Item(89, 13).Value = -12967.3335  # M89: -13590.25 -> -12967.3335
$ws.Cells.Item(89, 14).Value = -25689  # N89: -28673.25 -> -25689
$ws.Cells.Item(91, 8).Value = 9987.4  # H91: 9992 -> 9987.4
$ws.Cells.Item(91, 10).Value = 9987.4  # J91: 9992 -> 9987.4
$ws.Cells.Item(91, 12).Value = 9987.4  # L91: 9992 -> 9987.4
$ws.Cells.Item(91, 14).Value = -12795.4  # N91: -12800 -> -12795.4
$ws.Cells.Item(92, 8).Value = 526.25  # H92: 454.75 -> 526.25
$ws.Cells.Item(92, 9).Value = 387.14285  # I92: 295.8 -> 387.14285
$ws.Cells.Item(92, 10).Value = 1500  # J92: 1249.5 -> 1500
$ws.Cells.Item(92, 11).Value = 387.14285  # K92: 295.8 -> 387.14285
$ws.Cells.Item(92, 12).Value = 1500  # L92: 1249.5 -> 1500
$ws.Cells.Item(92, 13).Value = 860.85715  # M92: 952.2 -> 860.85715
$ws.Cells.Item(92, 14).Value = -3996  # N92: -3745.5 -> -3996
$ws.Cells.Item(101, 8).Value = 673.5714  # H101: 675.26666 -> 673.5714
$ws.Cells.Item(101, 9).Value = 781.875  # I101: 772.6667 -> 781.875
$ws.Cells.Item(101, 11).Value = 2345.625  # K101: 2318.0001 -> 2345.625
$ws.Cells.Item(101, 13).Value = -723.625  # M101: -696.0001000000002 -> -723.625
$ws.Cells.Item(107, 8).Value = 508.4375  # H107: 508.5 -> 508.4375
$ws.Cells.Item(107, 9).Value = 539.53845  # I107: 539.61536 -> 539.53845
$ws.Cells.Item(107, 11).Value = 539.53845  # K107: 539.61536 -> 539.53845
$ws.Cells.Item(107, 13).Value = 1380.46155  # M107: 1380.38464 -> 1380.46155
$ws.Cells.Item(112, 8).Value = 2633.7273  # H112: 2464.2307 -> 2633.7273
$ws.Cells.Item(112, 9).Value = 1295  # I112: 1290 -> 1295
$ws.Cells.Item(112, 10).Value = 2767.6  # J112: 2677.7273 -> 2767.6
$ws.Cells.Item(112, 11).Value = 3885  # K112: 3870 -> 3885
$ws.Cells.Item(112, 12).Value = 8302.799999999999  # L112: 8033.1819 -> 8302.799999999999
$ws.Cells.Item(112, 13).Value = -2777  # M112: -2762 -> -2777
$ws.Cells.Item(112, 14).Value = -10518.8  # N112: -10249.1819 -> -10518.8
$ws.Cells.Item(113, 8).Value = 9434.75  # H113: 9562.25 -> 9434.75
$ws.Cells.Item(113, 10).Value = 9746.5  # J113: 9916.5 -> 9746.5
$ws.Cells.Item(113, 12).Value = 9746.5  # L113: 9916.5 -> 9746.5
$ws.Cells.Item(113, 14).Value = -16254.5  # N113: -16424.5 -> -16254.5
$ws.Cells.Item(116, 8).Value = 4564.143  # H116: 4104.6 -> 4564.143
$ws.Cells.Item(116, 9).Value = 4909  # I116: 4506.375 -> 4909
$ws.Cells.Item(116, 10).Value = 2495  # J116: 2497.5 -> 2495
$ws.Cells.Item(116, 11).Value = 4909  # K116: 4506.375 -> 4909
$ws.Cells.Item(116, 12).Value = 2495  # L116: 2497.5 -> 2495
$ws.Cells.Item(116, 13).Value = -1467  # M116: -1064.375 -> -1467
$ws.Cells.Item(116, 14).Value = -9379  # N116: -9381.5 -> -9379
$ws.Cells.Item(132, 8).Value = 27029950  # H132: 27029994 -> 27029950
$ws.Cells.Item(132, 9).Value = 33336588  # I132: 34486148 -> 33336588
$ws.Cells.Item(132, 10).Value = 1499.2858  # J132: 1436.625 -> 1499.2858
$ws.Cells.Item(132, 11).Value = 100009764  # K132: 103458444 -> 100009764
$ws.Cells.Item(132, 12).Value = 4497.857400000001  # L132: 4309.875 -> 4497.857400000001
$ws.Cells.Item(132, 13).Value = -100007234  # M132: -103455914 -> -100007234
$ws.Cells.Item(132, 14).Value = -9557.857400000001  # N132: -9369.875 -> -9557.857400000001
$ws.Cells.Item(135, 8).Value = 677.8182  # H135: 661.0833 -> 677.8182
$ws.Cells.Item(135, 9).Value = 677.8182  # I135: 661.0833 -> 677.8182
$ws.Cells.Item(135, 11).Value = 6100.3638  # K135: 5949.7497 -> 6100.3638
$ws.Cells.Item(135, 13).Value = -3565.3638  # M135: -3414.7497 -> -3565.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5503.6523  # H32: 6381.275 -> 5503.6523
$ws.Cells.Item(32, 9).Value = 4130.6045  # I32: 4856.7295 -> 4130.6045
$ws.Cells.Item(32, 11).Value = 4130.6045  # K32: 4856.7295 -> 4130.6045
$ws.Cells.Item(32, 13).Value = -3843.6045  # M32: -4569.7295 -> -3843.6045
$ws.Cells.Item(74, 8).Value = 4602482  # H74: 3835901.2 -> 4602482
$ws.Cells.Item(74, 9).Value = 7668373  # I74: 6274669 -> 7668373
$ws.Cells.Item(74, 10).Value = 3645.6667  # J74: 3552.4285 -> 3645.6667
$ws.Cells.Item(74, 11).Value = 7668373  # K74: 6274669 -> 7668373
$ws.Cells.Item(74, 12).Value = 3645.6667  # L74: 3552.4285 -> 3645.6667
$ws.Cells.Item(74, 13).Value = -7667499  # M74: -6273795 -> -7667499
$ws.Cells.Item(74, 14).Value = -5393.6667  # N74: -5300.4285 -> -5393.6667
$ws.Cells.Item(77, 8).Value = 4602482  # H77: 3835901.2 -> 4602482
$ws.Cells.Item(77, 9).Value = 7668373  # I77: 6274669 -> 7668373
$ws.Cells.Item(77, 10).Value = 3645.6667  # J77: 3552.4285 -> 3645.6667
$ws.Cells.Item(77, 11).Value = 38341865  # K77: 31373345 -> 38341865
$ws.Cells.Item(77, 12).Value = 18228.3335  # L77: 17762.1425 -> 18228.3335
$ws.Cells.Item(77, 13).Value = -38337497  # M77: -31368977 -> -38337497
$ws.Cells.Item(77, 14).Value = -26964.3335  # N77: -26498.1425 -> -26964.3335
$ws.Cells.Item(97, 8).Value = 817.3333  # H97: 857.25 -> 817.3333
$ws.Cells.Item(97, 9).Value = 767.2105  # I97: 770.8421 -> 767.2105
$ws.Cells.Item(97, 10).Value = 1293.5  # J97: 2499 -> 1293.5
$ws.Cells.Item(97, 11).Value = 767.2105  # K97: 770.8421 -> 767.2105
$ws.Cells.Item(97, 12).Value = 1293.5  # L97: 2499 -> 1293.5
$ws.Cells.Item(97, 13).Value = -271.2105  # M97: -274.8421 -> -271.2105
$ws.Cells.Item(97, 14).Value = -2285.5  # N97: -3491 -> -2285.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 590.6667  # H22: 559 -> 590.6667
$ws.Cells.Item(22, 10).Value = 766  # J22: 800 -> 766
$ws.Cells.Item(22, 12).Value = 766  # L22: 800 -> 766
$ws.Cells.Item(22, 14).Value = -1112  # N22: -1146 -> -1112
$ws.Cells.Item(86, 8).Value = 13362942  # H86: 12376446 -> 13362942
$ws.Cells.Item(86, 9).Value = 31371.357  # I86: 33106.562 -> 31371.357
$ws.Cells.Item(86, 11).Value = 31371.357  # K86: 33106.562 -> 31371.357
$ws.Cells.Item(86, 13).Value = -30248.357  # M86: -31983.562 -> -30248.357
$ws.Cells.Item(89, 8).Value = 13362942  # H89: 12376446 -> 13362942
$ws.Cells.Item(89, 9).Value = 31371.357  # I89: 33106.562 -> 31371.357
$ws.Cells.Item(89, 11).Value = 156856.785  # K89: 165532.81 -> 156856.785
$ws.Cells.Item(89, 13).Value = -151240.785  # M89: -159916.81 -> -151240.785
$ws.Cells.Item(94, 8).Value = 1618.7  # H94: 1649.625 -> 1618.7
$ws.Cells.Item(94, 9).Value = 1648.5  # I94: 1699.6666 -> 1648.5
$ws.Cells.Item(94, 11).Value = 1648.5  # K94: 1699.6666 -> 1648.5
$ws.Cells.Item(94, 13).Value = -1197.5  # M94: -1248.6666 -> -1197.5
$ws.Cells.Item(99, 8).Value = 1699.75  # H99: 1933.3334 -> 1699.75
$ws.Cells.Item(99, 9).Value = 1699.75  # I99: 1933.3334 -> 1699.75
$ws.Cells.Item(99, 11).Value = 1699.75  # K99: 1933.3334 -> 1699.75
$ws.Cells.Item(99, 13).Value = -201.75  # M99: -435.3334 -> -201.75
$ws.Cells.Item(108, 8).Value = 69999.164  # H108: 69999.125 -> 69999.164
$ws.Cells.Item(108, 10).Value = 69999.164  # J108: 69999.125 -> 69999.164
$ws.Cells.Item(108, 12).Value = 69999.164  # L108: 69999.125 -> 69999.164
$ws.Cells.Item(108, 14).Value = -77679.164  # N108: -77679.125 -> -77679.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 2208.1667  # H5: 1511.75 -> 2208.1667
$ws.Cells.Item(5, 9).Value = 2269.8  # I5: 2120 -> 2269.8
$ws.Cells.Item(5, 10).Value = 1900  # J5: 498 -> 1900
$ws.Cells.Item(5, 11).Value = 2269.8  # K5: 2120 -> 2269.8
$ws.Cells.Item(5, 12).Value = 1900  # L5: 498 -> 1900
$ws.Cells.Item(5, 13).Value = -2157.8  # M5: -2008 -> -2157.8
$ws.Cells.Item(5, 14).Value = -2124  # N5: -722 -> -2124
$ws.Cells.Item(7, 8).Value = 52631760  # H7: 52631744 -> 52631760
$ws.Cells.Item(7, 9).Value = 66666890  # I7: 66666864 -> 66666890
$ws.Cells.Item(7, 11).Value = 66666890  # K7: 66666864 -> 66666890
$ws.Cells.Item(7, 13).Value = -66666777  # M7: -66666751 -> -66666777
$ws.Cells.Item(58, 8).Value = 3155.3333  # H58: 3510.9167 -> 3155.3333
$ws.Cells.Item(58, 9).Value = 3259.2856  # I58: 3510.9167 -> 3259.2856
$ws.Cells.Item(58, 10).Value = 1700  # J58: 0 -> 1700
$ws.Cells.Item(58, 11).Value = 3259.2856  # K58: 3510.9167 -> 3259.2856
$ws.Cells.Item(58, 12).Value = 1700  # L58: 0 -> 1700
$ws.Cells.Item(58, 13).Value = -3056.2856  # M58: -3307.9167 -> -3056.2856
$ws.Cells.Item(58, 14).Value = -2106  # N58: None -> -2106
$ws.Cells.Item(105, 8).Value = 1525  # H105: 1206.4 -> 1525
$ws.Cells.Item(105, 9).Value = 1525  # I105: 1206.4 -> 1525
$ws.Cells.Item(105, 11).Value = 1525  # K105: 1206.4 -> 1525
$ws.Cells.Item(105, 13).Value = 222  # M105: 540.5999999999999 -> 222
$ws.Cells.Item(136, 8).Value = 3155.3333  # H136: 3510.9167 -> 3155.3333
$ws.Cells.Item(136, 9).Value = 3259.2856  # I136: 3510.9167 -> 3259.2856
$ws.Cells.Item(136, 10).Value = 1700  # J136: 0 -> 1700
$ws.Cells.Item(136, 11).Value = 9777.856800000001  # K136: 10532.7501 -> 9777.856800000001
$ws.Cells.Item(136, 12).Value = 5100  # L136: 0 -> 5100
$ws.Cells.Item(136, 13).Value = -7227.856800000001  # M136: -7982.750100000001 -> -7227.856800000001
$ws.Cells.Item(136, 14).Value = -10200  # N136: None -> -10200
$ws.Cells.Item(141, 8).Value = 397037.1  # H141: 122475.625 -> 397037.1
$ws.Cells.Item(141, 10).Value = 397037.1  # J141: 122475.625 -> 397037.1
$ws.Cells.Item(141, 12).Value = 397037.1  # L141: 122475.625 -> 397037.1
$ws.Cells.Item(141, 14).Value = -407397.1  # N141: -132835.625 -> -407397.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 359.1154  # H107: 369.52 -> 359.1154
$ws.Cells.Item(107, 10).Value = 337.9091  # J107: 361.8 -> 337.9091
$ws.Cells.Item(107, 12).Value = 1013.7273  # L107: 1085.4 -> 1013.7273
$ws.Cells.Item(107, 14).Value = -4853.7273  # N107: -4925.4 -> -4853.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2103  # H80: 2492.2856 -> 2103
$ws.Cells.Item(80, 9).Value = 1193  # I80: 0 -> 1193
$ws.Cells.Item(80, 10).Value = 2204.111  # J80: 2492.2856 -> 2204.111
$ws.Cells.Item(80, 11).Value = 1193  # K80: 0 -> 1193
$ws.Cells.Item(80, 12).Value = 2204.111  # L80: 2492.2856 -> 2204.111
$ws.Cells.Item(80, 13).Value = -195  # M80: None -> -195
$ws.Cells.Item(80, 14).Value = -4200.111  # N80: -4488.2856 -> -4200.111
$ws.Cells.Item(83, 8).Value = 2103  # H83: 2492.2856 -> 2103
$ws.Cells.Item(83, 9).Value = 1193  # I83: 0 -> 1193
$ws.Cells.Item(83, 10).Value = 2204.111  # J83: 2492.2856 -> 2204.111
$ws.Cells.Item(83, 11).Value = 5965  # K83: 0 -> 5965
$ws.Cells.Item(83, 12).Value = 11020.555  # L83: 12461.428 -> 11020.555
$ws.Cells.Item(83, 13).Value = -973  # M83: None -> -973
$ws.Cells.Item(83, 14).Value = -21004.555  # N83: -22445.428 -> -21004.555
$ws.Cells.Item(122, 8).Value = 10007  # H122: 0 -> 10007
$ws.Cells.Item(122, 9).Value = 10007  # I122: 0 -> 10007
$ws.Cells.Item(122, 11).Value = 30021  # K122: 0 -> 30021
$ws.Cells.Item(122, 13).Value = -27571  # M122: None -> -27571
$ws.Cells.Item(132, 8).Value = 47621620  # H132: 33335762 -> 47621620
$ws.Cells.Item(132, 10).Value = 83335290  # J132: 47621060 -> 83335290
$ws.Cells.Item(132, 12).Value = 250005870  # L132: 142863180 -> 250005870
$ws.Cells.Item(132, 14).Value = -250010930  # N132: -142868240 -> -250010930

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 38468310  # H40: 37043884 -> 38468310
$ws.Cells.Item(40, 9).Value = 58829480  # I40: 58829468 -> 58829480
$ws.Cells.Item(40, 10).Value = 8321.888999999999  # J40: 8389.5 -> 8321.888999999999
$ws.Cells.Item(40, 11).Value = 58829480  # K40: 58829468 -> 58829480
$ws.Cells.Item(40, 12).Value = 8321.888999999999  # L40: 8389.5 -> 8321.888999999999
$ws.Cells.Item(40, 13).Value = -58829344  # M40: -58829332 -> -58829344
$ws.Cells.Item(40, 14).Value = -8593.888999999999  # N40: -8661.5 -> -8593.888999999999
$ws.Cells.Item(61, 8).Value = 3966.8667  # H61: 3999.5334 -> 3966.8667
$ws.Cells.Item(61, 9).Value = 4076  # I61: 4111 -> 4076
$ws.Cells.Item(61, 11).Value = 4076  # K61: 4111 -> 4076
$ws.Cells.Item(61, 13).Value = -3874  # M61: -3909 -> -3874
$ws.Cells.Item(113, 8).Value = 3966.8667  # H113: 3999.5334 -> 3966.8667
$ws.Cells.Item(113, 9).Value = 4076  # I113: 4111 -> 4076
$ws.Cells.Item(113, 11).Value = 4076  # K113: 4111 -> 4076
$ws.Cells.Item(113, 13).Value = -1906  # M113: -1941 -> -1906
$ws.Cells.Item(140, 8).Value = 62323  # H140: 62323.5 -> 62323
$ws.Cells.Item(140, 10).Value = 62323  # J140: 62323.5 -> 62323
$ws.Cells.Item(140, 12).Value = 62323  # L140: 62323.5 -> 62323
$ws.Cells.Item(140, 14).Value = -72683  # N140: -72683.5 -> -72683

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 7202.067  # H81: 8070.387 -> 7202.067
$ws.Cells.Item(81, 9).Value = 3640.4736  # I81: 4009.8948 -> 3640.4736
$ws.Cells.Item(81, 10).Value = 13353.909  # J81: 14499.5 -> 13353.909
$ws.Cells.Item(81, 11).Value = 7280.9472  # K81: 8019.7896 -> 7280.9472
$ws.Cells.Item(81, 12).Value = 26707.818  # L81: 28999 -> 26707.818
$ws.Cells.Item(81, 13).Value = -6219.9472  # M81: -6958.7896 -> -6219.9472
$ws.Cells.Item(81, 14).Value = -28829.818  # N81: -31121 -> -28829.818
$ws.Cells.Item(84, 8).Value = 7202.067  # H84: 8070.387 -> 7202.067
$ws.Cells.Item(84, 9).Value = 3640.4736  # I84: 4009.8948 -> 3640.4736
$ws.Cells.Item(84, 10).Value = 13353.909  # J84: 14499.5 -> 13353.909
$ws.Cells.Item(84, 11).Value = 36404.736  # K84: 40098.948 -> 36404.736
$ws.Cells.Item(84, 12).Value = 133539.09  # L84: 144995 -> 133539.09
$ws.Cells.Item(84, 13).Value = -31100.736  # M84: -34794.948 -> -31100.736
$ws.Cells.Item(84, 14).Value = -144147.09  # N84: -155603 -> -144147.09
$ws.Cells.Item(113, 8).Value = 644.2308  # H113: 675.0769 -> 644.2308
$ws.Cells.Item(113, 9).Value = 623.5454999999999  # I113: 660 -> 623.5454999999999
$ws.Cells.Item(113, 11).Value = 1870.6365  # K113: 1980 -> 1870.6365
$ws.Cells.Item(113, 13).Value = 299.3635000000002  # M113: 190 -> 299.3635000000002
$ws.Cells.Item(119, 8).Value = 70000  # H119: 0 -> 70000
$ws.Cells.Item(119, 10).Value = 70000  # J119: 0 -> 70000
$ws.Cells.Item(119, 12).Value = 70000  # L119: 0 -> 70000
$ws.Cells.Item(119, 14).Value = -79676  # N119: None -> -79676
$ws.Cells.Item(122, 8).Value = 2638.1538  # H122: 2279.889 -> 2638.1538
$ws.Cells.Item(122, 9).Value = 2441.3333  # I122: 2216.8 -> 2441.3333
$ws.Cells.Item(122, 10).Value = 5000  # J122: 2595.3333 -> 5000
$ws.Cells.Item(122, 11).Value = 7323.999899999999  # K122: 6650.400000000001 -> 7323.999899999999
$ws.Cells.Item(122, 12).Value = 15000  # L122: 7785.999899999999 -> 15000
$ws.Cells.Item(122, 13).Value = -4873.999899999999  # M122: -4200.400000000001 -> -4873.999899999999
$ws.Cells.Item(122, 14).Value = -19900  # N122: -12685.9999 -> -19900
$ws.Cells.Item(126, 8).Value = 2439.4  # H126: 3474254.2 -> 2439.4
$ws.Cells.Item(126, 9).Value = 2316.1667  # I126: 5210032 -> 2316.1667
$ws.Cells.Item(126, 10).Value = 2624.25  # J126: 2699.25 -> 2624.25
$ws.Cells.Item(126, 11).Value = 6948.500100000001  # K126: 15630096 -> 6948.500100000001
$ws.Cells.Item(126, 12).Value = 7872.75  # L126: 8097.75 -> 7872.75
$ws.Cells.Item(126, 13).Value = -4478.500100000001  # M126: -15627626 -> -4478.500100000001
$ws.Cells.Item(126, 14).Value = -12812.75  # N126: -13037.75 -> -12812.75
